$wb = $excel.ActiveWorkbook

# The handoff status text changed: "Ready for handoff" -> "In Translation".
# That shared string is used on the Overview sheet (columns E "zh-cn" and
# F "de-de") and on the per-locale "zh-cn"/"de-de" sheets (column C "Status").
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2:F4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2:C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2:C4").Value = "In Translation"

# The status report re-fits the columns that held the status text, which
# shrink now that "In Translation" is shorter than "Ready for handoff".
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
